# Localize the instruction strings from English to Dutch and move the
# active selection to C2 (matching the author's re-save in Dutch Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Duid aan welke van de twee afbeeldingen jij verkiest"
$ws.Range("B2").Value = "De afbeelding links: "
$ws.Range("C2").Value = "de afbeelding rechts:"

[void]$ws.Range("C2").Select()
